$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 3 new item rows into the alphabetically-sorted table (rows 7-36),
# copying the format of an existing data row (row 7) so that merges, styles
# and borders match the rest of the table. Rows are inserted top-to-bottom,
# with later insertion points expressed in terms of the row numbers that
# result *after* the earlier inserts have already shifted things down:
#
#   - "AMIGRAINE ADCO 30 TABLETS"                     -> before old row 8  (BI-KETOGESIC)
#   - "CLEAREST 14 CAPS"                               -> before old row 13 (DECLOPHEN), which is row 14 after the first insert
#   - "FLIXONASE 50 MCG/METERED NASAL SPRAY 120 DOSE"  -> before old row 17 (GARAMYCIN), which is row 19 after the first two inserts
# ---------------------------------------------------------------------------

$ws.Rows("7:7").Copy()
$ws.Rows("8:8").Insert()

$ws.Rows("7:7").Copy()
$ws.Rows("14:14").Insert()

$ws.Rows("7:7").Copy()
$ws.Rows("19:19").Insert()

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Fill in the data for the three newly-inserted rows.
# Columns: C = item name, H = current balance, L = order limit,
#          N = price, P = sale price, Q = transaction count
# ---------------------------------------------------------------------------

$ws.Range("C8").Value = "AMIGRAINE ADCO 30 TABLETS"
$ws.Range("H8").Value = "13:2"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "93.00"
$ws.Range("P8").Value = "30.6900"
$ws.Range("Q8").Value = "0:1"

$ws.Range("C14").Value = "CLEAREST 14 CAPS"
$ws.Range("H14").Value = "1:1"
$ws.Range("L14").Value = "1"
$ws.Range("N14").Value = "44.00"
$ws.Range("P14").Value = "22.0000"
$ws.Range("Q14").Value = "0:1"

$ws.Range("C19").Value = "FLIXONASE 50 MCG/METERED NASAL SPRAY 120 DOSE"
$ws.Range("H19").Value = "1:0"
$ws.Range("L19").Value = "1"
$ws.Range("N19").Value = "133.00"
$ws.Range("P19").Value = "133.0000"
$ws.Range("Q19").Value = "1:0"

# ---------------------------------------------------------------------------
# Renumber column A (the running index "م") for the whole data block, which
# now spans rows 7-39 (33 items, was 30).
# ---------------------------------------------------------------------------

for ($i = 0; $i -lt 33; $i++) {
    $r = 7 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
}

# ---------------------------------------------------------------------------
# Update the grand-total cell (now on row 40, was row 37) and the generated
# timestamp in the footer (now on row 41, was row 38).
# ---------------------------------------------------------------------------

$ws.Range("P40").Value = 2314.25
$ws.Range("A41").Value = "Sunday, 25 May, 2025 1:07 PM"
